# Add season-record columns (Wins, Losses, Ties) to the team stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (AC1) onto the
# three new header cells so they match the existing bold/bordered look.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Season record values for every player row (2 through 53)
$lastRow = 53
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 97   # AD -> Wins
    $ws.Cells.Item($r, 31).Value2 = 65   # AE -> Losses
    $ws.Cells.Item($r, 32).Value2 = 0    # AF -> Ties
}
